# Edit: (1) swap the slide design's theme colours from "Integral" back to
# the default "Office Theme" palette, and (2) change the table style used
# by the cash-flow table on slide 16 from the "Medium"-style built-in id to
# the "No Style, No Grid" built-in id.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Theme: Integral -> Office Theme
#
# The deck's slide master currently carries the "Integral" theme colours.
# Re-apply the twelve standard "Office Theme" colours (dk1/lt1/dk2/lt2/
# accent1-6/hlink/folHlink, in that order) through the presentation's
# shared theme colour scheme so every slide (they all share one master)
# picks up the Office palette again.
$officeThemeColors = @(
    0x000000,  # dk1      - Black
    0xFFFFFF,  # lt1      - White
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $hex = $officeThemeColors[$i - 1]
    $r = ($hex -shr 16) -band 0xFF
    $g = ($hex -shr 8) -band 0xFF
    $b = $hex -band 0xFF
    $themeColors.Colors($i).RGB = $r -bor ($g -shl 8) -bor ($b -shl 16)
}

# ---------------------------------------------------------------------
# 2) Table style on slide 16's table: switch to "No Style, No Grid"
$slide16 = $p.Slides.Item(16)
for ($i = 1; $i -le $slide16.Shapes.Count; $i++) {
    $shp = $slide16.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{5D1CAA9B-E50A-4A47-994F-F466D23302EC}")
    }
}
